# Add team record (Wins/Losses/Ties) columns to the roster sheet.
# Mirrors xml_diff: new header cells AD1:AF1 ("Wins","Losses","Ties")
# styled like the other header cells, and AD:AF filled with 97/64/0
# for every data row (2-44).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (AC1, which
# carries the bold/centered/bordered header style) onto the three new
# header cells so they pick up the same style index instead of Excel
# creating a brand-new one.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every data row (2 through 44) gets the same team record.
$ws.Range("AD2:AD44").Value = 97
$ws.Range("AE2:AE44").Value = 64
$ws.Range("AF2:AF44").Value = 0
